# Results of discussion 16-4
# Add a new "Organisation aggregate" class above "Industry" on the
# Classes sheet, and re-parent "Industry" under it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Classes")

# Insert a new row above the current "Industry" row (row 10), copying
# formatting down from the row above (row 9, "Organisation").
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the "Organisation aggregate" class.
$ws.Cells.Item(10, 2).Value = "Oragnisation aggregate"
$ws.Cells.Item(10, 3).Value = "An aggregate of organisations."
$ws.Cells.Item(10, 4).Value = "Object aggregate"

# "Industry" (now row 11) is re-parented under "Organisation aggregate".
$ws.Cells.Item(11, 4).Value = "Organisation aggregate"
